$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 82
$ws.Range("A2").Value = 84
$ws.Range("A3").ClearContents()
